$wb = $excel.ActiveWorkbook

# --- Test_Case_3: Analytical solution comparison (GWpath-style, plain numbers) ---
$ws3 = $wb.Worksheets.Item("Test_Case_3")
$ws3.Range("B2").Value = 1282.669543993601
$ws3.Range("C2").Value = 0.5972768771205283
$ws3.Range("F2").Value = 3832.8930403376
$ws3.Range("G2").Value = 5.601162800581742
$ws3.Range("I2").Value = -3828.310450015801
$ws3.Range("J2").Value = 5.720698050606194

# --- Test_Case_4: Modpath6 vs Percent Difference table (numeric-looking text cells) ---
$ws4 = $wb.Worksheets.Item("Test_Case_4")
$ws4.Range("B3:G4").NumberFormat = "@"
$ws4.Range("B3").Value = "9354317.212192766"
$ws4.Range("C3").Value = "214.74706024030934"
$ws4.Range("D3").Value = "2355.3772450506603"
$ws4.Range("E3").Value = "2197.4601353337093"
$ws4.Range("F3").Value = "5941.83402009447"
$ws4.Range("G3").Value = "5976.18170204034"
$ws4.Range("B4").Value = "7.71726067689379"
$ws4.Range("C4").Value = "7.71726067689379"
$ws4.Range("D4").Value = "2.9135969025881248"
$ws4.Range("E4").Value = "5.908756275426588"
$ws4.Range("F4").Value = "0.5976046774849153"
$ws4.Range("G4").Value = "4.254768786046714"
$ws4.Range("B3:G4").Style = "Normal"

# --- Test_Case_5: Modpath6 vs Percent Difference table (numeric-looking text cells) ---
$ws5 = $wb.Worksheets.Item("Test_Case_5")
$ws5.Range("B3:G4").NumberFormat = "@"
$ws5.Range("B3").Value = "9367480.367901502"
$ws5.Range("C3").Value = "215.0492468059148"
$ws5.Range("D3").Value = "2354.5218369839713"
$ws5.Range("E3").Value = "2038.3388591967523"
$ws5.Range("F3").Value = "5941.090629167855"
$ws5.Range("G3").Value = "5817.6381606943905"
$ws5.Range("B4").Value = "7.647366367730347"
$ws5.Range("C4").Value = "7.647366367730347"
$ws5.Range("D4").Value = "2.8954503050796774"
$ws5.Range("E4").Value = "2.156961705001133"
$ws5.Range("F4").Value = "0.6038604110002143"
$ws5.Range("G4").Value = "5.595861539126026"
$ws5.Range("B3:G4").Style = "Normal"
